$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 59.7
$ws.Range("I2").Value = 55.22222
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 55.22222
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 57.77778
$ws.Range("N2").Value = -326

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 831.375
$ws.Range("I18").Value = 831.375
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 831.375
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -547.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 19608836
$ws.Range("I32").Value = 47619588
$ws.Range("J32").Value = 1307.6
$ws.Range("K32").Value = 47619588
$ws.Range("L32").Value = 1307.6
$ws.Range("M32").Value = -47619262
$ws.Range("N32").Value = -1959.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1129.5714
$ws.Range("I53").Value = 61.4
$ws.Range("J53").Value = 3800
$ws.Range("K53").Value = 61.4
$ws.Range("L53").Value = 3800
$ws.Range("M53").Value = 575.6
$ws.Range("N53").Value = -5074

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 828.5
$ws.Range("I129").Value = 439.375
$ws.Range("J129").Value = 902.619
$ws.Range("K129").Value = 1318.125
$ws.Range("L129").Value = 2707.857
$ws.Range("M129").Value = 3681.875
$ws.Range("N129").Value = -12707.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1792493.4
$ws.Range("I137").Value = 4808674.5
$ws.Range("J137").Value = 5126.7036
$ws.Range("K137").Value = 14426023.5
$ws.Range("L137").Value = 15380.1108
$ws.Range("M137").Value = -14423473.5
$ws.Range("N137").Value = -20480.1108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1919.4642
$ws.Range("I2").Value = 1919.4642
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1919.4642
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1806.4642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 4388334.5
$ws.Range("I6").Value = 6537502
$ws.Range("J6").Value = 90000
$ws.Range("K6").Value = 6537502
$ws.Range("L6").Value = 90000
$ws.Range("M6").Value = -6537329
$ws.Range("N6").Value = -90346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1582.909
$ws.Range("I45").Value = 1581
$ws.Range("J45").Value = 1584.5
$ws.Range("K45").Value = 1581
$ws.Range("L45").Value = 1584.5
$ws.Range("M45").Value = -1204
$ws.Range("N45").Value = -2338.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1904.8536
$ws.Range("I74").Value = 1624.258
$ws.Range("J74").Value = 2774.7
$ws.Range("K74").Value = 1624.258
$ws.Range("L74").Value = 2774.7
$ws.Range("M74").Value = -750.258
$ws.Range("N74").Value = -4522.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1904.8536
$ws.Range("I77").Value = 1624.258
$ws.Range("J77").Value = 2774.7
$ws.Range("K77").Value = 8121.29
$ws.Range("L77").Value = 13873.5
$ws.Range("M77").Value = -3753.29
$ws.Range("N77").Value = -22609.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1919.4642
$ws.Range("I116").Value = 1919.4642
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1919.4642
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 374.5358000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1919.4642
$ws.Range("I3").Value = 1919.4642
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1919.4642
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1805.4642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 29933.334
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 29933.334
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 29933.334
$ws.Range("N81").Value = -32055.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 29933.334
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 29933.334
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 89800.00199999999
$ws.Range("N84").Value = -100408.002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1906.1765
$ws.Range("I107").Value = 1599.4286
$ws.Range("J107").Value = 3337.6667
$ws.Range("K107").Value = 1599.4286
$ws.Range("L107").Value = 3337.6667
$ws.Range("M107").Value = 320.5714
$ws.Range("N107").Value = -7177.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 28973.8
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 28973.8
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 28973.8
$ws.Range("N41").Value = -29829.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 100031090
$ws.Range("I51").Value = 500000000
$ws.Range("J51").Value = 38859.75
$ws.Range("K51").Value = 500000000
$ws.Range("L51").Value = 38859.75
$ws.Range("M51").Value = -499999264
$ws.Range("N51").Value = -40331.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 30412.8
$ws.Range("I59").Value = 40104
$ws.Range("J59").Value = 27990
$ws.Range("K59").Value = 40104
$ws.Range("L59").Value = 27990
$ws.Range("M59").Value = -38959
$ws.Range("N59").Value = -30280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 10847.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 100031090
$ws.Range("I61").Value = 500000000
$ws.Range("J61").Value = 38859.75
$ws.Range("K61").Value = 500000000
$ws.Range("L61").Value = 38859.75
$ws.Range("M61").Value = -499999652
$ws.Range("N61").Value = -39555.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 564.56757
$ws.Range("I107").Value = 433.7857
$ws.Range("J107").Value = 971.44446
$ws.Range("K107").Value = 433.7857
$ws.Range("L107").Value = 971.44446
$ws.Range("M107").Value = 1486.2143
$ws.Range("N107").Value = -4811.44446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6517.5293
$ws.Range("I5").Value = 6799.875
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 20399.625
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -20287.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1165.6522
$ws.Range("I86").Value = 900.25
$ws.Range("J86").Value = 1307.2
$ws.Range("K86").Value = 2700.75
$ws.Range("L86").Value = 3921.6
$ws.Range("M86").Value = -1514.75
$ws.Range("N86").Value = -6293.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1165.6522
$ws.Range("I89").Value = 900.25
$ws.Range("J89").Value = 1307.2
$ws.Range("K89").Value = 8102.25
$ws.Range("L89").Value = 11764.8
$ws.Range("M89").Value = -2174.25
$ws.Range("N89").Value = -23620.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10001.091
$ws.Range("I122").Value = 300
$ws.Range("J122").Value = 10971.2
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 98740.8
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -103640.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6517.5293
$ws.Range("I135").Value = 6799.875
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 61198.875
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -58663.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1429.2307
$ws.Range("I113").Value = 1711.25
$ws.Range("J113").Value = 978
$ws.Range("K113").Value = 1711.25
$ws.Range("L113").Value = 978
$ws.Range("M113").Value = 458.75
$ws.Range("N113").Value = -5318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 24287.705
$ws.Range("I126").Value = 48649.125
$ws.Range("J126").Value = 2633.111
$ws.Range("K126").Value = 145947.375
$ws.Range("L126").Value = 7899.333
$ws.Range("M126").Value = -143477.375
$ws.Range("N126").Value = -12839.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2308.9524
$ws.Range("I7").Value = 1817.5264
$ws.Range("J7").Value = 6977.5
$ws.Range("K7").Value = 1817.5264
$ws.Range("L7").Value = 6977.5
$ws.Range("M7").Value = -1705.5264

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3630.2104
$ws.Range("I40").Value = 2666.6
$ws.Range("J40").Value = 7243.75
$ws.Range("K40").Value = 2666.6
$ws.Range("L40").Value = 7243.75
$ws.Range("M40").Value = -2530.6
$ws.Range("N40").Value = -7515.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2308.9524
$ws.Range("I126").Value = 1817.5264
$ws.Range("J126").Value = 6977.5
$ws.Range("K126").Value = 5452.5792
$ws.Range("L126").Value = 20932.5
$ws.Range("M126").Value = -2982.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1875.96
$ws.Range("I136").Value = 1415.15
$ws.Range("J136").Value = 3719.2
$ws.Range("K136").Value = 4245.450000000001
$ws.Range("L136").Value = 11157.6
$ws.Range("M136").Value = -1695.450000000001
$ws.Range("N136").Value = -16257.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 28410
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 28410
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 28410
$ws.Range("N75").Value = -30282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 28410
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 28410
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 85230
$ws.Range("N78").Value = -94590
